$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45958
$ws.Range("B2").Value = 108.28
$ws.Range("C2").Value = 105.77
$ws.Range("D2").Value = 105.5
$ws.Range("E2").Value = 104.28
$ws.Range("F2").Value = 104.29
$ws.Range("G2").Value = 99.62
$ws.Range("H2").Value = 110.21
$ws.Range("I2").Value = 121.49
$ws.Range("J2").Value = 111.83
$ws.Range("K2").Value = 100.44
$ws.Range("L2").Value = 75.2
$ws.Range("M2").Value = 62.03
$ws.Range("N2").Value = 60
$ws.Range("O2").Value = 63.5
$ws.Range("P2").Value = 70.81
$ws.Range("Q2").Value = 74.97
$ws.Range("R2").Value = 98.86
$ws.Range("S2").Value = 110.92
$ws.Range("T2").Value = 128.04
$ws.Range("U2").Value = 130.13
$ws.Range("V2").Value = 142.24
$ws.Range("W2").Value = 130.04
$ws.Range("X2").Value = 111.33
$ws.Range("Y2").Value = 106.7
$ws.Range("Z2").Value = 101.52
$ws.Range("AB2").Value = 122.58
$ws.Range("AD2").Value = 136.14
$ws.Range("AF2").Value = 129.08
$ws.Range("AG2").Value = "5h-16h"
